$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 3 first
$ws.Range("B3").Value = "Adhiyan"
$ws.Range("C3").Value = "School Kid"

# Update row 2: replace "rajappa " with "Vanaja"
$ws.Range("B2").Value = "Vanaja"

# Column B width adjustment (auto-fit after data change widened column B)
$ws.Columns.Item(2).ColumnWidth = 6.83

# Update selection to D3
$ws.Range("D3").Select()
